# "Add box to baseline"
#
# The "Interventions" sheet lists intervention rows twice: once "applied to"
# the Baseline (Calibration) run and once "applied to" the Hypothetical
# Scenario run. In this edit the Social Distancing entries are pulled out of
# the Baseline block entirely (Baseline keeps only Handwashing x2,
# Self-isolation, Screening and Household Isolation, with refreshed coverage
# %), while the Hypothetical Scenario block keeps every intervention plus its
# own brand-new "Social Distancing" box appended at the end (new dates /
# coverage). The two rows this vacates (15-16) go blank, and the unused
# trailing blank rows (47-50) are trimmed from the bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Interventions")

# --- Final contents for rows 2-14 (Intervention, Start, End, Coverage %, Apply to) ---
$rows = @(
    @(2,  "Handwashing",                    "2/15/2020", "4/15/2020",  30, "Baseline (Calibration)"),
    @(3,  "Handwashing",                    "4/16/2020", "12/31/2020", 60, "Baseline (Calibration)"),
    @(4,  "Self-isolation if Symptomatic",   "2/15/2020", "3/31/2020",  60, "Baseline (Calibration)"),
    @(5,  "Screening (when S.I.)",           "2/15/2020", "3/31/2020",  60, "Baseline (Calibration)"),
    @(6,  "Household Isolation (when S.I.)", "2/15/2020", "3/31/2020",  50, "Baseline (Calibration)"),
    @(7,  "Handwashing",                    "2/15/2020", "4/15/2020",  30, "Hypothetical Scenario"),
    @(8,  "Handwashing",                    "4/16/2020", "12/31/2020", 60, "Hypothetical Scenario"),
    @(9,  "Self-isolation if Symptomatic",   "2/15/2020", "3/31/2020",  60, "Hypothetical Scenario"),
    @(10, "Screening (when S.I.)",           "2/15/2020", "3/31/2020",  60, "Hypothetical Scenario"),
    @(11, "Household Isolation (when S.I.)", "2/15/2020", "3/31/2020",  50, "Hypothetical Scenario"),
    @(12, "International Travel Ban",        "2/15/2020", "12/31/2020", 100, "Hypothetical Scenario"),
    @(13, "Social Distancing",               "2/5/2020",  "2/29/2020",  15, "Hypothetical Scenario"),
    @(14, "Social Distancing",               "3/1/2020",  "6/30/2020",  30, "Hypothetical Scenario")
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Range("A$r").Value = $row[1]
    $ws.Range("B$r").Value = $row[2]
    $ws.Range("C$r").Value = $row[3]
    $ws.Range("D$r").Value = $row[4]
    $ws.Range("E$r").Value = $row[5]
}

# --- Rows 15 & 16 no longer carry an intervention (blank them out, keep styles) ---
$ws.Range("A15:E16").ClearContents() | Out-Null

# --- Trim the unused trailing blank rows; table now ends at row 46 ---
$ws.Range("A47:E50").Delete() | Out-Null

# --- Interventions becomes the active sheet / tab with E15 selected ---
$ws.Activate() | Out-Null
$ws.Range("E15").Select() | Out-Null
